$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.988.16"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.211.21"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.56%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "2.547.36"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.832"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").Value = "2.207.93"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "41.866.78"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000106"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.21%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0792"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -12.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "64.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.196"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "2.421.12"
$ws.Range("E51").Value = "  -2.10%  "
